$wb = $excel.ActiveWorkbook

# --- ALC (sheet index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H43").Value = 1429.5294
$ws.Range("I43").Value = 300
$ws.Range("J43").Value = 1500.125
$ws.Range("K43").Value = 300
$ws.Range("L43").Value = 1500.125
$ws.Range("M43").Value = -231
$ws.Range("N43").Value = -1638.125

# --- ARM (sheet index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H109").Value = 29000
$ws.Range("J109").Value = 29000
$ws.Range("L109").Value = 29000
$ws.Range("N109").Value = -31774

$ws.Range("H132").Value = 41143.137
$ws.Range("I132").Value = 28604.703
$ws.Range("J132").Value = 74280.42999999999
$ws.Range("K132").Value = 85814.109
$ws.Range("L132").Value = 222841.29
$ws.Range("M132").Value = -83284.109
$ws.Range("N132").Value = -227901.29

# --- BSM (sheet index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 1684
$ws.Range("I94").Value = 1684
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1684
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1233
$ws.Range("N94").ClearContents()

# --- CRP (sheet index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H4").Value = 14100
$ws.Range("J4").Value = 14100
$ws.Range("L4").Value = 14100
$ws.Range("N4").Value = -14324

$ws.Range("H31").Value = 1884.26
$ws.Range("I31").Value = 961.15216
$ws.Range("J31").Value = 12500
$ws.Range("K31").Value = 961.15216
$ws.Range("L31").Value = 12500
$ws.Range("M31").Value = -666.15216
$ws.Range("N31").Value = -13090

$ws.Range("H34").Value = 1884.26
$ws.Range("I34").Value = 961.15216
$ws.Range("J34").Value = 12500
$ws.Range("K34").Value = 961.15216
$ws.Range("L34").Value = 12500
$ws.Range("M34").Value = -759.15216
$ws.Range("N34").Value = -12904

$ws.Range("H62").Value = 2796.5833
$ws.Range("I62").Value = 2796.5833
$ws.Range("K62").Value = 2796.5833
$ws.Range("M62").Value = -2172.5833

$ws.Range("H63").Value = 30270.666
$ws.Range("J63").Value = 30270.666
$ws.Range("L63").Value = 30270.666
$ws.Range("N63").Value = -31642.666

$ws.Range("H65").Value = 2796.5833
$ws.Range("I65").Value = 2796.5833
$ws.Range("K65").Value = 13982.9165
$ws.Range("M65").Value = -10862.9165

$ws.Range("H66").Value = 30270.666
$ws.Range("J66").Value = 30270.666
$ws.Range("L66").Value = 90811.99800000001
$ws.Range("N66").Value = -97675.99800000001

# --- CUL (sheet index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H4").Value = 9092150
$ws.Range("J4").Value = 20001600
$ws.Range("L4").Value = 60004800
$ws.Range("N4").Value = -60005024

$ws.Range("H62").Value = 4166.6665
$ws.Range("J62").Value = 4166.6665
$ws.Range("L62").Value = 12499.9995
$ws.Range("N62").Value = -13871.9995

$ws.Range("H63").Value = 3575
$ws.Range("I63").Value = 1000
$ws.Range("J63").Value = 4433.3335
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 13300.0005
$ws.Range("M63").Value = -2251
$ws.Range("N63").Value = -14798.0005

$ws.Range("H64").Value = 3716.1428
$ws.Range("I64").Value = 950
$ws.Range("J64").Value = 3928.923
$ws.Range("K64").Value = 2850
$ws.Range("L64").Value = 11786.769
$ws.Range("M64").Value = -2580
$ws.Range("N64").Value = -12326.769

$ws.Range("H65").Value = 4166.6665
$ws.Range("J65").Value = 4166.6665
$ws.Range("L65").Value = 37499.9985
$ws.Range("N65").Value = -44363.9985

$ws.Range("H66").Value = 3575
$ws.Range("I66").Value = 1000
$ws.Range("J66").Value = 4433.3335
$ws.Range("K66").Value = 9000
$ws.Range("L66").Value = 39900.0015
$ws.Range("M66").Value = -5256
$ws.Range("N66").Value = -47388.0015

$ws.Range("H67").Value = 3716.1428
$ws.Range("I67").Value = 950
$ws.Range("J67").Value = 3928.923
$ws.Range("K67").Value = 2850
$ws.Range("L67").Value = 11786.769
$ws.Range("M67").Value = -1914
$ws.Range("N67").Value = -13658.769

$ws.Range("H75").Value = 2473.2856
$ws.Range("I75").Value = 2437.3333
$ws.Range("J75").Value = 2500.25
$ws.Range("K75").Value = 7311.999899999999
$ws.Range("L75").Value = 7500.75
$ws.Range("M75").Value = -6313.999899999999
$ws.Range("N75").Value = -9496.75

$ws.Range("H76").Value = 2809.2307
$ws.Range("I76").Value = 1640
$ws.Range("J76").Value = 3160
$ws.Range("K76").Value = 4920
$ws.Range("L76").Value = 9480
$ws.Range("M76").Value = -4537
$ws.Range("N76").Value = -10246

$ws.Range("H78").Value = 2473.2856
$ws.Range("I78").Value = 2437.3333
$ws.Range("J78").Value = 2500.25
$ws.Range("K78").Value = 21935.9997
$ws.Range("L78").Value = 22502.25
$ws.Range("M78").Value = -16943.9997
$ws.Range("N78").Value = -32486.25

$ws.Range("H79").Value = 2809.2307
$ws.Range("I79").Value = 1640
$ws.Range("J79").Value = 3160
$ws.Range("K79").Value = 4920
$ws.Range("L79").Value = 9480
$ws.Range("M79").Value = -3594
$ws.Range("N79").Value = -12132

$ws.Range("H131").Value = 738.5333000000001
$ws.Range("J131").Value = 996.3333
$ws.Range("L131").Value = 2988.9999
$ws.Range("N131").Value = -13068.9999

# --- GSM (sheet index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 2954.5925
$ws.Range("I80").Value = 2316.7646
$ws.Range("J80").Value = 4038.9
$ws.Range("K80").Value = 2316.7646
$ws.Range("L80").Value = 4038.9
$ws.Range("M80").Value = -1318.7646
$ws.Range("N80").Value = -6034.9

$ws.Range("H83").Value = 2954.5925
$ws.Range("I83").Value = 2316.7646
$ws.Range("J83").Value = 4038.9
$ws.Range("K83").Value = 11583.823
$ws.Range("L83").Value = 20194.5
$ws.Range("M83").Value = -6591.823
$ws.Range("N83").Value = -30178.5

$ws.Range("H93").Value = 22251
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 22251
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 22251
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -25995

$ws.Range("H139").Value = 82000
$ws.Range("J139").Value = 82000
$ws.Range("L139").Value = 82000
$ws.Range("N139").Value = -92280

# --- LTW (sheet index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 3185.7144
$ws.Range("I61").Value = 3185.7144
$ws.Range("K61").Value = 3185.7144
$ws.Range("M61").Value = -2983.7144

$ws.Range("H82").Value = 2525
$ws.Range("I82").Value = 2550
$ws.Range("J82").Value = 2500
$ws.Range("K82").Value = 2550
$ws.Range("L82").Value = 2500
$ws.Range("M82").Value = -2189
$ws.Range("N82").Value = -3222

$ws.Range("H85").Value = 2525
$ws.Range("I85").Value = 2550
$ws.Range("J85").Value = 2500
$ws.Range("K85").Value = 2550
$ws.Range("L85").Value = 2500
$ws.Range("M85").Value = -1302
$ws.Range("N85").Value = -4996

$ws.Range("H93").Value = 1342.9286
$ws.Range("I93").Value = 1009.3
$ws.Range("K93").Value = 1009.3
$ws.Range("M93").Value = 238.7

$ws.Range("H95").Value = 60000
$ws.Range("J95").Value = 60000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -65492

$ws.Range("H113").Value = 3185.7144
$ws.Range("I113").Value = 3185.7144
$ws.Range("K113").Value = 3185.7144
$ws.Range("M113").Value = -1015.7144

# --- WVR (sheet index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41802

$ws.Range("H97").Value = 60000
$ws.Range("J97").Value = 60000
$ws.Range("L97").Value = 60000
$ws.Range("N97").Value = -61982
